$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price/volume figures to the latest scrape; keep the numeric-looking
# price strings stored as text (matching the workbook's existing inlineStr convention)
# instead of letting Excel auto-convert them into real numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D26", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "265.20"
$ws.Range("D3").Value = "22.81"
$ws.Range("D4").Value = "6.230"
$ws.Range("D5").Value = "0.06171"
$ws.Range("D6").Value = "3.563"
$ws.Range("D7").Value = "6.704"
$ws.Range("D8").Value = "1.359"
$ws.Range("D9").Value = "0.8166"
$ws.Range("D11").Value = "0.08194"
$ws.Range("D13").Value = "0.03149"
$ws.Range("D14").Value = "0.09262"
$ws.Range("D15").Value = "3.891"
$ws.Range("D16").Value = "0.001686"
$ws.Range("D17").Value = "0.04840"
$ws.Range("D18").Value = "0.0006262"
$ws.Range("D19").Value = "0.006169"
$ws.Range("D20").Value = "0.006270"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("D23").Value = "3.699"
$ws.Range("D26").Value = "0.1247"
$ws.Range("D27").Value = "0.0002682"
$ws.Range("D40").Value = "0.04588"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.007218"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1135"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.003226"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "0.01043"
$ws.Range("D45").Value = "0.00006137"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").Value = "0.7703"
$ws.Range("D48").Value = "0.1947"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").Value = "0.01240"
